$d = $word.ActiveDocument

# Asher's line after the "Must've been rough" beat: swap his expression
# tag from "smiling_nervous" to "smiling_eyes_closed".
$d.Content.Find.Execute(
    "Asher (neutral smiling_nervous): Must" + [char]8217 + "ve been rough.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Asher (neutral smiling_eyes_closed): Must" + [char]8217 + "ve been rough.",
    2
)

# Asher's "Yikes." line: swap his expression tag from "thoughtful" to
# "smiling_nervous".
$d.Content.Find.Execute(
    "Asher (neutral thoughtful) : Yikes.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Asher (neutral smiling_nervous) : Yikes.",
    2
)
